$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 2.429935333333333
$ws.Range("N2").Value = 7.289806
$ws.Range("O2").Value = 0.0007999299918632063
$ws.Range("P2").Value = 0.0007999299918632063
$ws.Range("Q2").Value = 279.9804764981164
$ws.Range("R2").Value = 2519.824288483048
$ws.Range("S2").Value = 0.0002229731854343324
$ws.Range("T2").Value = 0.0002229731854343324

$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.001089498967837074
$ws.Range("P3").Value = 0.001089498967837074
$ws.Range("Q3").Value = 381.3314205768569
$ws.Range("R3").Value = 3431.982785191712
$ws.Range("S3").Value = 0.0003036878950121828
$ws.Range("T3").Value = 0.0003036878950121827

$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 1.214141
$ws.Range("N4").Value = 3.642423
$ws.Range("O4").Value = 0.0003996928588706414
$ws.Range("P4").Value = 0.0003996928588706414
$ws.Range("Q4").Value = 139.8949885837426
$ws.Range("R4").Value = 1259.054897253684
$ws.Range("S4").Value = 0.0001114107369948223
$ws.Range("T4").Value = 0.0001114107369948222

$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 3030.731364
$ws.Range("N5").Value = 9092.194092
$ws.Range("O5").Value = 0.9977108781814291
$ws.Range("P5").Value = 0.9977108781814292
$ws.Range("Q5").Value = 349205.0178415611
$ws.Range("R5").Value = 3142845.16057405
$ws.Range("S5").Value = 0.2781028026370602
$ws.Range("T5").Value = 0.2781028026370602

$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 2.429935333333333
$ws.Range("N6").Value = 7.289806
$ws.Range("O6").Value = 0.0007999299918632063
$ws.Range("P6").Value = 0.0007999299918632063
$ws.Range("Q6").Value = 449.153625756804
$ws.Range("R6").Value = 4042.382631811236
$ws.Range("S6").Value = 0.0003577007080529357
$ws.Range("T6").Value = 0.0003577007080529357

$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.001089498967837074
$ws.Range("P7").Value = 0.001089498967837074
$ws.Range("S7").Value = 0.0004871858239875921
$ws.Range("T7").Value = 0.0004871858239875921

$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 1.214141
$ws.Range("N8").Value = 3.642423
$ws.Range("O8").Value = 0.0003996928588706414
$ws.Range("P8").Value = 0.0003996928588706414
$ws.Range("Q8").Value = 224.424010322082
$ws.Range("R8").Value = 2019.816092898738
$ws.Range("S8").Value = 0.0001787286638530982
$ws.Range("T8").Value = 0.0001787286638530982

$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 3030.731364
$ws.Range("N9").Value = 9092.194092
$ws.Range("O9").Value = 0.9977108781814291
$ws.Range("P9").Value = 0.9977108781814292
$ws.Range("Q9").Value = 560205.8466996779
$ws.Range("R9").Value = 5041852.620297101
$ws.Range("S9").Value = 0.4461414013573364
$ws.Range("T9").Value = 0.4461414013573364

$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 2.429935333333333
$ws.Range("N10").Value = 7.289806
$ws.Range("O10").Value = 0.0007999299918632063
$ws.Range("P10").Value = 0.0007999299918632063
$ws.Range("Q10").Value = 147.1349090714689
$ws.Range("R10").Value = 1324.21418164322
$ws.Range("S10").Value = 0.0001171765252155965
$ws.Range("T10").Value = 0.0001171765252155965

$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.001089498967837074
$ws.Range("P11").Value = 0.001089498967837074
$ws.Range("Q11").Value = 200.3967012072978
$ws.Range("R11").Value = 1803.57031086568
$ws.Range("S11").Value = 0.0001595935951592107
$ws.Range("T11").Value = 0.0001595935951592107

$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 1.214141
$ws.Range("N12").Value = 3.642423
$ws.Range("O12").Value = 0.0003996928588706414
$ws.Range("P12").Value = 0.0003996928588706414
$ws.Range("Q12").Value = 73.51739907822333
$ws.Range("R12").Value = 661.65659170401
$ws.Range("S12").Value = 0.000058548399025347
$ws.Range("T12").Value = 0.00005854839902534699

$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 3030.731364
$ws.Range("N13").Value = 9092.194092
$ws.Range("O13").Value = 0.9977108781814291
$ws.Range("P13").Value = 0.9977108781814292
$ws.Range("Q13").Value = 183513.6834898716
$ws.Range("R13").Value = 1651623.151408844
$ws.Range("S13").Value = 0.1461481567940677
$ws.Range("T13").Value = 0.1461481567940677

$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 2.429935333333333
$ws.Range("N14").Value = 7.289806
$ws.Range("O14").Value = 0.0007999299918632063
$ws.Range("P14").Value = 0.0007999299918632063
$ws.Range("Q14").Value = 128.1781371086611
$ws.Range("R14").Value = 1153.60323397795
$ws.Range("S14").Value = 0.0001020795731603416
$ws.Range("T14").Value = 0.0001020795731603416

$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.001089498967837074
$ws.Range("P15").Value = 0.001089498967837074
$ws.Range("Q15").Value = 174.5777124244222
$ws.Range("R15").Value = 1571.1994118198
$ws.Range("S15").Value = 0.0001390316536780882
$ws.Range("T15").Value = 0.0001390316536780882

$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 1.214141
$ws.Range("N16").Value = 3.642423
$ws.Range("O16").Value = 0.0003996928588706414
$ws.Range("P16").Value = 0.0003996928588706414
$ws.Range("Q16").Value = 64.04546221144166
$ws.Range("R16").Value = 576.409159902975
$ws.Range("S16").Value = 0.00005100505899737398
$ws.Range("T16").Value = 0.00005100505899737398

$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 3030.731364
$ws.Range("N17").Value = 9092.194092
$ws.Range("O17").Value = 0.9977108781814291
$ws.Range("P17").Value = 0.9977108781814292
$ws.Range("Q17").Value = 159869.8924145491
$ws.Range("R17").Value = 1438829.031730942
$ws.Range("S17").Value = 0.1273185173929648
$ws.Range("T17").Value = 0.1273185173929648
